$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '23.192.06'
Set-TextValue "E2" '  +0.36%  '
Set-TextValue "D3" '1.601.47'
Set-TextValue "E4" '  +0.03%  '
Set-TextValue "E5" '  +0.04%  '
Set-TextValue "D6" '303.35'
Set-TextValue "E6" '  +0.51%  '
Set-TextValue "E7" '  +0.01%  '
Set-TextValue "D8" '52.05'
Set-TextValue "E8" '  +4.62%  '
Set-TextValue "D9" '0.3614'
Set-TextValue "E9" '  -1.03%  '
Set-TextValue "D10" '1.268'
Set-TextValue "E10" '  -0.06%  '
Set-TextValue "E11" '  +0.01%  '
Set-TextValue "D12" '0.08116'
Set-TextValue "D13" '22.71'
Set-TextValue "E13" '  -1.80%  '
Set-TextValue "D14" '6.574'
Set-TextValue "E14" '  -0.30%  '
Set-TextValue "D15" '7.401'
Set-TextValue "E15" '  -0.12%  '
Set-TextValue "D16" '0.00001243'
Set-TextValue "E16" '  -1.43%  '
Set-TextValue "D17" '1.601.72'
Set-TextValue "E17" '  +0.14%  '
Set-TextValue "D18" '94.14'
Set-TextValue "E18" '  +2.61%  '
Set-TextValue "D19" '0.06877'
Set-TextValue "E19" '  +0.52%  '
Set-TextValue "E20" '  -2.18%  '
Set-TextValue "D21" '6.539'
Set-TextValue "E21" '  -0.74%  '
Set-TextValue "D22" '1.001'
Set-TextValue "E22" '  +0.01%  '
Set-TextValue "E23" '  -0.44%  '
Set-TextValue "D24" '23.197.46'
Set-TextValue "E24" '  +0.36%  '
Set-TextValue "D25" '2.400'
Set-TextValue "E25" '  +2.56%  '
Set-TextValue "D26" '2.988'
Set-TextValue "E26" '  +9.57%  '
Set-TextValue "E27" '  +0.29%  '
Set-TextValue "D28" '149.36'
Set-TextValue "E28" '  -0.52%  '
Set-TextValue "D29" '5.247'
Set-TextValue "E29" '  -0.46%  '
Set-TextValue "D30" '134.06'
Set-TextValue "E30" '  +1.11%  '
Set-TextValue "D31" '2.385'
Set-TextValue "E31" '  -0.21%  '
Set-TextValue "D32" '6.766'
Set-TextValue "E32" '  -1.53%  '
Set-TextValue "D33" '1.779.22'
Set-TextValue "E33" '  +0.25%  '
Set-TextValue "D34" '0.9687'
Set-TextValue "E34" '  +0.70%  '
Set-TextValue "D35" '0.07503'
Set-TextValue "E35" '  -2.72%  '
Set-TextValue "D36" '10.27'
Set-TextValue "E36" '  +1.79%  '
Set-TextValue "E37" '  -0.69%  '
Set-TextValue "E38" '  -2.25%  '
Set-TextValue "D39" '0.08802'
Set-TextValue "E39" '  -1.03%  '
Set-TextValue "D40" '6.089'
Set-TextValue "E40" '  -3.26%  '
Set-TextValue "D41" '0.7100'
Set-TextValue "E41" '  +0.12%  '
Set-TextValue "D42" '1.360'
Set-TextValue "E42" '  -0.77%  '
Set-TextValue "D43" '12.48'
Set-TextValue "E43" '  -1.20%  '
Set-TextValue "D44" '15.65'
Set-TextValue "E44" '  +1.78%  '
Set-TextValue "D45" '0.6521'
Set-TextValue "E45" '  -1.46%  '
Set-TextValue "D46" '2.311'
Set-TextValue "E46" '  -0.13%  '
Set-TextValue "D47" '4.019'
Set-TextValue "E47" '  +0.63%  '
Set-TextValue "D48" '132.07'
Set-TextValue "E48" '  +0.13%  '
Set-TextValue "D49" '0.07968'
Set-TextValue "E49" '  +0.42%  '
Set-TextValue "E50" '  -2.11%  '
Set-TextValue "D51" '1.211'
Set-TextValue "E51" '  +1.17%  '
